$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows 10-12: numeric score updates (Right/Wrong/NotAttempt/Max + Total) ---
$ws.Range("B10").Value2 = 19
$ws.Range("C10").Value2 = 2
$ws.Range("D10").Value2 = 7
$ws.Range("E10").Value2 = 28
$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1
$ws.Range("B12").Value2 = 76
$ws.Range("C12").Value2 = -2

# E12: "Absent" -> computed score fraction
$ws.Range("E12").Value2 = "74/112"

# A10/A11/A12 labels ("No." / "Marking" / "Total") get the bold header style
# (copy format from A9, which already carries that style); text itself is unchanged
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null

# --- Column A (rows 16-40): show the student answer, colored green if correct or red if wrong ---
# "Correct" cells copy their format from B10 (an existing green "correctStyle" cell)
$ws.Range("B10").Copy() | Out-Null
$greenCellsA = @("A16", "A17", "A18", "A22", "A23", "A25", "A29", "A30", "A31", "A32", "A33", "A36", "A37", "A38", "A39", "A40")
foreach ($c in $greenCellsA) { $ws.Range($c).PasteSpecial(-4122) | Out-Null }

# "Incorrect" cells copy their format from C10 (an existing red "incorrectStyle" cell)
$ws.Range("C10").Copy() | Out-Null
$redCellsA = @("A19", "A26")
foreach ($c in $redCellsA) { $ws.Range($c).PasteSpecial(-4122) | Out-Null }

# Now write the actual answer text into the colored A cells
$ws.Range("A16").Value2 = "Option A"
$ws.Range("A17").Value2 = "Option D"
$ws.Range("A18").Value2 = "Option B"
$ws.Range("A19").Value2 = "Option D"
$ws.Range("A22").Value2 = "Option D"
$ws.Range("A23").Value2 = "Option D"
$ws.Range("A25").Value2 = "Option A"
$ws.Range("A26").Value2 = "Option B"
$ws.Range("A29").Value2 = "Option D"
$ws.Range("A30").Value2 = "Option B"
$ws.Range("A31").Value2 = "Option D"
$ws.Range("A32").Value2 = "Option C"
$ws.Range("A33").Value2 = "Option D"
$ws.Range("A36").Value2 = "Option A"
$ws.Range("A37").Value2 = "Option A"
$ws.Range("A38").Value2 = "Option A"
$ws.Range("A39").Value2 = "Option D"
$ws.Range("A40").Value2 = "Option D"

# --- Column D (rows 16-18): new second-question-block answers, all correct (green) ---
$ws.Range("B10").Copy() | Out-Null
$greenCellsD = @("D16", "D17", "D18")
foreach ($c in $greenCellsD) { $ws.Range($c).PasteSpecial(-4122) | Out-Null }
$ws.Range("D16").Value2 = "Option A"
$ws.Range("D17").Value2 = "Option C"
$ws.Range("D18").Value2 = "Option D"

# --- Drop the now-unused cells: D/E below row 18, and the whole third block G/H ---
$ws.Range("D19:E40").Clear() | Out-Null
$ws.Range("G15:H40").Clear() | Out-Null

Write-Output "edit applied"
